$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 7, 14, 20, 23, 30, 33, 39, 43, 45, 46, 48, 55, 58, 62, 64, 71, 72, 73, 75, 81, 83, 86, 91, 93, 94, 96, 97, 98, 99, 105, 106, 112, 113, 114, 118, 119, 120, 124, 131, 132, 134, 137, 139, 140, 144, 150, 153, 158, 159, 160, 164, 167, 170, 171, 173, 174, 175, 177, 179, 181, 182, 187, 189)

foreach ($r in $rows) {
    $ws.Range("G${r}:N${r}").Value = "N/A"
}
